$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source data
# (many of the price strings look numeric, e.g. "1.008", "0.06636", etc.)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "20.511.80"
$ws.Range("E2").Value = "  -7.07%  "

# Row 3
$ws.Range("D3").Value = "1.451.32"
$ws.Range("E3").Value = "  -6.95%  "

# Row 4
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.86%  "

# Row 5
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6
$ws.Range("D6").Value = "277.75"
$ws.Range("E6").Value = "  -4.85%  "

# Row 7
$ws.Range("D7").Value = "0.3713"
$ws.Range("E7").Value = "  -6.52%  "

# Row 8
$ws.Range("D8").Value = "0.3095"
$ws.Range("E8").Value = "  -4.60%  "

# Row 9
$ws.Range("D9").Value = "41.43"
$ws.Range("E9").Value = "  -6.64%  "

# Row 10
$ws.Range("D10").Value = "1.016"
$ws.Range("E10").Value = "  -6.21%  "

# Row 11
$ws.Range("D11").Value = "0.06636"
$ws.Range("E11").Value = "  -8.89%  "

# Row 12
$ws.Range("D12").Value = "1.009"
$ws.Range("E12").Value = "  +0.88%  "

# Row 13
$ws.Range("D13").Value = "5.436"
$ws.Range("E13").Value = "  -4.97%  "

# Row 14
$ws.Range("D14").Value = "17.39"
$ws.Range("E14").Value = "  -7.85%  "

# Row 15
$ws.Range("D15").Value = "6.182"
$ws.Range("E15").Value = "  -7.25%  "

# Row 16
$ws.Range("D16").Value = "1.454.85"
$ws.Range("E16").Value = "  -7.03%  "

# Row 17
$ws.Range("D17").Value = "0.00001025"
$ws.Range("E17").Value = "  -8.89%  "

# Row 18
$ws.Range("D18").Value = "0.06077"
$ws.Range("E18").Value = "  -7.88%  "

# Row 19
$ws.Range("D19").Value = "77.88"
$ws.Range("E19").Value = "  -7.22%  "

# Row 20
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").Value = "  +0.81%  "

# Row 21
$ws.Range("D21").Value = "5.750"
$ws.Range("E21").Value = "  -8.66%  "

# Row 22
$ws.Range("D22").Value = "14.66"
$ws.Range("E22").Value = "  -6.17%  "

# Row 23
$ws.Range("D23").Value = "11.00"
$ws.Range("E23").Value = "  -3.18%  "

# Row 24
$ws.Range("D24").Value = "2.307"
$ws.Range("E24").Value = "  -2.37%  "

# Row 25
$ws.Range("D25").Value = "20.517.03"
$ws.Range("E25").Value = "  -7.09%  "

# Row 26
$ws.Range("D26").Value = "2.263"
$ws.Range("E26").Value = "  -6.88%  "

# Row 27
$ws.Range("D27").Value = "143.33"
$ws.Range("E27").Value = "  -3.46%  "

# Row 28
$ws.Range("D28").Value = "17.22"
$ws.Range("E28").Value = "  -7.85%  "

# Row 29
$ws.Range("D29").Value = "1.619.33"
$ws.Range("E29").Value = "  -6.93%  "

# Row 30
$ws.Range("D30").Value = "109.66"
$ws.Range("E30").Value = "  -8.32%  "

# Row 31
$ws.Range("B31").Value = "HuobiToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D31").Value = "3.710"
$ws.Range("E31").Value = "  -23.92%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.9172"
$ws.Range("E32").Value = "  -7.43%  "

# Row 33
$ws.Range("D33").Value = "5.475"
$ws.Range("E33").Value = "  -7.95%  "

# Row 34
$ws.Range("D34").Value = "0.07768"
$ws.Range("E34").Value = "  -6.89%  "

# Row 35
$ws.Range("D35").Value = "8.314"
$ws.Range("E35").Value = "  -9.48%  "

# Row 36
$ws.Range("D36").Value = "1.439"
$ws.Range("E36").Value = "  -10.60%  "

# Row 37
$ws.Range("D37").Value = "11.05"
$ws.Range("E37").Value = "  +2.25%  "

# Row 38
$ws.Range("D38").Value = "1.008"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39
$ws.Range("D39").Value = "4.788"
$ws.Range("E39").Value = "  -7.24%  "

# Row 40
$ws.Range("D40").Value = "0.05648"
$ws.Range("E40").Value = "  -6.43%  "

# Row 41
$ws.Range("D41").Value = "0.02055"
$ws.Range("E41").Value = "  -10.02%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.130"
$ws.Range("E42").Value = "  -6.73%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1924"
$ws.Range("E43").Value = "  -6.71%  "

# Row 44
$ws.Range("D44").Value = "3.594"
$ws.Range("E44").Value = "  -4.66%  "

# Row 45
$ws.Range("D45").Value = "0.5357"
$ws.Range("E45").Value = "  -8.39%  "

# Row 46
$ws.Range("D46").Value = "12.18"
$ws.Range("E46").Value = "  -7.28%  "

# Row 47
$ws.Range("D47").Value = "0.5185"
$ws.Range("E47").Value = "  -7.67%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "111.06"
$ws.Range("E48").Value = "  -6.36%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.786"
$ws.Range("E49").Value = "  -6.21%  "

# Row 50
$ws.Range("D50").Value = "1.065"
$ws.Range("E50").Value = "  -6.77%  "

# Row 51
$ws.Range("D51").Value = "0.06340"
$ws.Range("E51").Value = "  -7.13%  "
